$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.123.93"
$ws.Range("E2").Value = "  +0.54%  "
$ws.Range("D3").Value = "2.519.36"
$ws.Range("E3").Value = "  +0.88%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'535.79"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.11%  "
$ws.Range("D6").Value = "'139.63"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.04%  "
$ws.Range("E7").Value = "  +0.24%  "
$ws.Range("E8").Value = "  -1.02%  "
$ws.Range("D9").Value = "2.520.19"
$ws.Range("E9").Value = "  -0.16%  "
$ws.Range("D10").Value = "'0.0993"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.18%  "
$ws.Range("E11").Value = "  +1.56%  "
$ws.Range("D12").Value = "'5.42"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.89%  "
$ws.Range("E13").Value = "  +1.66%  "
$ws.Range("D14").Value = "2.963.50"
$ws.Range("E14").Value = "  +0.91%  "
$ws.Range("D15").Value = "59.066.86"
$ws.Range("E15").Value = "  +0.55%  "
$ws.Range("D16").Value = "'22.95"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.75%  "
$ws.Range("E17").Value = "  +1.52%  "
$ws.Range("D18").Value = "2.519.37"
$ws.Range("E18").Value = "  +0.36%  "
$ws.Range("D19").Value = "'10.90"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.27%  "
$ws.Range("D20").Value = "'4.23"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.21%  "
$ws.Range("D21").Value = "'321.52"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.18%  "
$ws.Range("E22").Value = "  -0.16%  "
$ws.Range("D23").Value = "'5.82"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.68%  "
$ws.Range("D24").Value = "'62.87"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.95%  "
$ws.Range("D25").Value = "'0.424"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.20%  "
$ws.Range("E26").Value = "  +1.71%  "
$ws.Range("D27").Value = "'0.999"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.55%  "
$ws.Range("E28").Value = "  +1.12%  "
$ws.Range("D29").Value = "'6.72"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.25%  "
$ws.Range("B30").Value = "PEPE"
$ws.Range("C30").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D30").Value = "0.0₃0767"
$ws.Range("E30").Value = "  +0.53%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "'1.80"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.89%  "
$ws.Range("D32").Value = "'160.58"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.25%  "
$ws.Range("E33").Value = "  +0.27%  "
$ws.Range("E34").Value = "  -3.81%  "
$ws.Range("D35").Value = "'1.46"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.83%  "
$ws.Range("E36").Value = "  -0.24%  "
$ws.Range("D37").Value = "'4.20"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.31%  "
$ws.Range("E38").Value = "  -1.42%  "
$ws.Range("D39").Value = "'36.96"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.18%  "
$ws.Range("E40").Value = "  +0.32%  "
$ws.Range("E41").Value = "  -0.04%  "
$ws.Range("D42").Value = "'5.24"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.92%  "
$ws.Range("D43").Value = "'282.92"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.69%  "
$ws.Range("E44").Value = "  +0.29%  "
$ws.Range("E45").Value = "  +0.96%  "
$ws.Range("E46").Value = "  -1.17%  "
$ws.Range("D47").Value = "'0.0930"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.61%  "
$ws.Range("D48").Value = "'122.82"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.42%  "
$ws.Range("D49").Value = "'18.50"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.02%  "
$ws.Range("E50").Value = "  -0.08%  "
$ws.Range("E51").Value = "  -1.70%  "
